$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.996.25"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.256.87"
$ws.Range("E3").Value = "  -3.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.Value = "'213.33"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'628.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.54%  "
$c = $ws.Range("D7")
$c.Value = "'0.388"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +12.38%  "
$c = $ws.Range("D8")
$c.Value = "'0.710"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +14.75%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.251.76"
$ws.Range("E10").Value = "  -3.37%  "
$c = $ws.Range("D11")
$c.Value = "'0.578"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.70%  "
$c = $ws.Range("D12")
$c.Value = "'0.188"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +12.37%  "
$c = $ws.Range("D13")
$c.Value = "'0.0000269"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D14")
$c.Value = "'5.51"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D15")
$c.Value = "'34.36"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").Value = "3.858.03"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("D17").Value = "87.729.85"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "3.282.19"
$ws.Range("E18").Value = "  -2.61%  "
$c = $ws.Range("D19")
$c.Value = "'3.25"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "
$c = $ws.Range("D20")
$c.Value = "'14.08"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.66%  "
$c = $ws.Range("D21")
$c.Value = "'437.60"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -7.43%  "
$c = $ws.Range("D22")
$c.Value = "'9.02"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "
$c = $ws.Range("D23")
$c.Value = "'5.35"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.77%  "
$c = $ws.Range("D24")
$c.Value = "'7.41"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$c = $ws.Range("D25")
$c.Value = "'5.38"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.21%  "
$c = $ws.Range("D26")
$c.Value = "'12.50"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -8.72%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D27")
$c.Value = "'0.0000143"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +11.29%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.436.16"
$ws.Range("E28").Value = "  -2.04%  "
$c = $ws.Range("D29")
$c.Value = "'77.39"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.82%  "
$c = $ws.Range("D30")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "
$c = $ws.Range("D31")
$c.Value = "'0.177"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -14.05%  "
$ws.Range("E32").Value = "  +0.28%  "
$c = $ws.Range("D33")
$c.Value = "'8.91"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -5.15%  "
$c = $ws.Range("D34")
$c.Value = "'571.31"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.23%  "
$c = $ws.Range("D35")
$c.Value = "'7.32"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("E36").Value = "  -10.51%  "
$ws.Range("E37").Value = "  -5.14%  "
$c = $ws.Range("D38")
$c.Value = "'0.140"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -8.69%  "
$c = $ws.Range("D39")
$c.Value = "'23.02"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.63%  "
$c = $ws.Range("D40")
$c.Value = "'3.31"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.88%  "
$ws.Range("E41").Value = "  +0.43%  "
$c = $ws.Range("D42")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$c = $ws.Range("D43")
$c.Value = "'0.405"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.58%  "
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("E45").Value = "  -0.07%  "
$c = $ws.Range("D46")
$c.Value = "'151.70"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.62%  "
$c = $ws.Range("D47")
$c.Value = "'0.136"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +19.59%  "
$c = $ws.Range("D48")
$c.Value = "'180.67"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.68%  "
$c = $ws.Range("D49")
$c.Value = "'45.23"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("E50").Value = "  -3.01%  "
$c = $ws.Range("D51")
$c.Value = "'4.27"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.03%  "
